$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix mojibake degree sign in J5 / J6 (UTF-8 bytes mis-decoded as Latin-1 -> correct to U+00B0) ---
$degree = [char]0x00B0
$ws.Range("J5").Value = $degree
$ws.Range("J6").Value = $degree

# --- Widen column J (10th column) from 6 to 7 characters (raw OOXML width) ---
# ColumnWidth 5.17 <-> raw width 6 ; ColumnWidth 6.17 <-> raw width 7
$ws.Columns.Item(10).ColumnWidth = 6.17

# --- Add new blank separator row 40 (present but empty, like row 11) ---
# Touch a harmless row-level property so the engine materializes the row element
# without adding any cell or attribute (matches "<row r="40"/>" in the target).
$ws.Rows.Item(40).OutlineLevel = 0

# --- Add new message block rows 41-45 : "DV_DRIVING_DYNAMICS_2" (ID 0x501) ---

# Row 41: message title/id header -> copy formatting+values from the matching header row 1
$ws.Range("A1:B1").Copy($ws.Range("A41:B41"))
$ws.Range("A41").Value = "Message: DV_DRIVING_DYNAMICS_2"
$ws.Range("B41").Value = "ID: 0x501"

# Row 42: column titles -> copy formatting+values from the matching column-title row 2
$ws.Range("A2:K2").Copy($ws.Range("A42:K42"))

# Rows 43-45: signal data -> copy formatting (and blank-template values) from row 9,
# which has the same "empty Min/Max, empty Choices, signed" shape we need here.
$ws.Range("A9:K9").Copy($ws.Range("A43:K43"))
$ws.Range("A9:K9").Copy($ws.Range("A44:K44"))
$ws.Range("A9:K9").Copy($ws.Range("A45:K45"))

# Row 43: ACCELERATION_LONGITUDINAL
$ws.Range("A43").Value = "ACCELERATION_LONGITUDINAL"
$ws.Range("B43").Value = 0
$ws.Range("C43").Value = 16
$ws.Range("D43").Value = "Intel"
$ws.Range("E43").Value = $true
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 0
$ws.Range("H43").ClearContents()
$ws.Range("I43").ClearContents()
$ws.Range("J43").Value = "m/s^2"
$ws.Range("K43").ClearContents()

# Row 44: ACCELERATION_LATERAL
$ws.Range("A44").Value = "ACCELERATION_LATERAL"
$ws.Range("B44").Value = 16
$ws.Range("C44").Value = 16
$ws.Range("D44").Value = "Intel"
$ws.Range("E44").Value = $true
$ws.Range("F44").Value = 1
$ws.Range("G44").Value = 0
$ws.Range("H44").ClearContents()
$ws.Range("I44").ClearContents()
$ws.Range("J44").Value = "m/s^2"
$ws.Range("K44").ClearContents()

# Row 45: YAW_RATE  (unit text contains a literal embedded line break: "m/s" + newline + "2")
$ws.Range("A45").Value = "YAW_RATE"
$ws.Range("B45").Value = 32
$ws.Range("C45").Value = 16
$ws.Range("D45").Value = "Intel"
$ws.Range("E45").Value = $true
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 0
$ws.Range("H45").ClearContents()
$ws.Range("I45").ClearContents()
$ws.Range("J45").Value = "m/s`n2"
$ws.Range("K45").ClearContents()
